$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period (fiscal year) column headers ---
# Shift the 5-year window forward by one year (drop 1396/12, add 1401/12)
$ws.Cells.Item(8, 4).Value = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 5).Value = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 8).Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ---
$ws.Cells.Item(9, 4).Value = "1399-04-19 (8)"
$ws.Cells.Item(9, 5).Value = "1400-04-14 (8)"
$ws.Cells.Item(9, 6).Value = "1401-04-15 (10)"
$ws.Cells.Item(9, 7).Value = "1402-02-27 (9)"
$ws.Cells.Item(9, 8).Value = "1402-02-27 (2)"

# --- Data rows: reset the read_price derived figures (D:H) to 0 ---
# Rows whose five values are all plain zeroes
$zeroRows = 11,12,13,14,16,17,19,20,22,24,25,26,27
foreach ($r in $zeroRows) {
    for ($c = 4; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# Row 15: D=0, E:H = "-"
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = "-"
$ws.Cells.Item(15, 6).Value = "-"
$ws.Cells.Item(15, 7).Value = "-"
$ws.Cells.Item(15, 8).Value = "-"

# Row 18: D:H all = "-"
$ws.Cells.Item(18, 4).Value = "-"
$ws.Cells.Item(18, 5).Value = "-"
$ws.Cells.Item(18, 6).Value = "-"
$ws.Cells.Item(18, 7).Value = "-"
$ws.Cells.Item(18, 8).Value = "-"

# Row 21: D="-", E:H = 0
$ws.Cells.Item(21, 4).Value = "-"
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0

# Row 23: D:H all = "-"
$ws.Cells.Item(23, 4).Value = "-"
$ws.Cells.Item(23, 5).Value = "-"
$ws.Cells.Item(23, 6).Value = "-"
$ws.Cells.Item(23, 7).Value = "-"
$ws.Cells.Item(23, 8).Value = "-"
